# TC09_C3DC_phs002431_TrtmntType-Imunothrpy.xlsx
# Commit message: "Updated remaining queries for C3DC"
#
# The author finished renaming the join columns used by every saved SQL
# query on Sheet1 (StatQuery in C2, plus the per-tab TabQuery cells B2:B7):
#   std.id / prt.id / "study.id" / "participant.id"
#     -> std.study_id / prt.participant_id / "study.study_id" / "participant.participant_id"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Fix-Query($text) {
    $text = $text.Replace('std.id = prt."study.id"', 'std.study_id = prt."study.study_id"')
    $text = $text.Replace('prt.id = dgn."participant.id"', 'prt.participant_id = dgn."participant.participant_id"')
    $text = $text.Replace('prt.id = trt."participant.id"', 'prt.participant_id = trt."participant.participant_id"')
    $text = $text.Replace('prt.id = trr."participant.id"', 'prt.participant_id = trr."participant.participant_id"')
    $text = $text.Replace('prt.id = srv."participant.id"', 'prt.participant_id = srv."participant.participant_id"')
    $text = $text.Replace('std.id = rfs."study.id"', 'std.study_id = rfs."study.study_id"')
    return $text
}

$cells = @("C2", "B2", "B3", "B4", "B5", "B6", "B7")
foreach ($addr in $cells) {
    $r = $ws.Range($addr)
    $r.Value2 = Fix-Query $r.Value2
}

# Column C (the empty "paste query output here" column) was widened after the
# edit, from a best-fit 60.83 chars to a fixed 69.33 chars.
$ws.Columns.Item(3).ColumnWidth = 68.5

# The author ended the session with the TreatmentRespTab query (B6) visible at
# the top of the window and C7 (SurvivalTab's result column) selected.
$excel.ActiveWindow.ScrollRow = 6
$ws.Range("C7").Select()
